# Updated cryptos list data (Price / Volume(1h) columns) for Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '58.703.29'
$ws.Range("E2").Value = '  -2.70%  '

# Row 3
$ws.Range("D3").Value = '2.656.09'
$ws.Range("E3").Value = '  -0.99%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.35'
$ws.Range("E5").Value = '  -0.11%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.88'
$ws.Range("E6").Value = '  -1.87%  '

# Row 7
$ws.Range("E7").Value = '  +0.22%  '

# Row 8
$ws.Range("E8").Value = '  -1.31%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.90'
$ws.Range("E9").Value = '  +6.79%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.102'
$ws.Range("E10").Value = '  -3.66%  '

# Row 11
$ws.Range("E11").Value = '  -1.33%  '

# Row 13
$ws.Range("D13").Value = '3.124.06'
$ws.Range("E13").Value = '  -0.94%  '

# Row 14
$ws.Range("D14").Value = '58.702.38'
$ws.Range("E14").Value = '  -2.68%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.96'
$ws.Range("E15").Value = '  -1.98%  '

# Row 16
$ws.Range("E16").Value = '  -1.82%  '

# Row 17
$ws.Range("D17").Value = '2.644.67'
$ws.Range("E17").Value = '  -1.58%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '338.53'
$ws.Range("E18").Value = '  -3.86%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.38'
$ws.Range("E19").Value = '  -3.54%  '

# Row 20
$ws.Range("E20").Value = '  -1.55%  '

# Row 21
$ws.Range("E21").Value = '  +1.36%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.12%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.92'
$ws.Range("E23").Value = '  +1.18%  '

# Row 24
$ws.Range("E24").Value = '  +0.37%  '

# Row 25
$ws.Range("E25").Value = '  -1.74%  '

# Row 26
$ws.Range("E26").Value = '  +0.88%  '

# Row 27
$ws.Range("D27").Value = '0.0₃0803'
$ws.Range("E27").Value = '  -1.91%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.13'
$ws.Range("E28").Value = '  -3.12%  '

# Row 29
$ws.Range("E29").Value = '  -1.58%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.06%  '

# Row 31
$ws.Range("E31").Value = '  -0.13%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.87'
$ws.Range("E32").Value = '  -1.64%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '150.99'
$ws.Range("E33").Value = '  +2.49%  '

# Row 34
$ws.Range("E34").Value = '  -3.18%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.935'
$ws.Range("E35").Value = '  -2.21%  '

# Row 36
$ws.Range("E36").Value = '  -6.01%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.876'
$ws.Range("E37").Value = '  -0.35%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.79'
$ws.Range("E38").Value = '  -0.38%  '

# Row 39
$ws.Range("E39").Value = '  -5.22%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.58'
$ws.Range("E40").Value = '  -3.30%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.38%  '

# Row 42
$ws.Range("E42").Value = '  -0.31%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '275.85'
$ws.Range("E43").Value = '  -3.08%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.68'
$ws.Range("E44").Value = '  -1.99%  '

# Row 45
$ws.Range("E45").Value = '  -2.63%  '

# Row 46
$ws.Range("E46").Value = '  +2.03%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0536'
$ws.Range("E47").Value = '  -1.07%  '

# Row 48
$ws.Range("D48").Value = '2.059.80'
$ws.Range("E48").Value = '  -3.55%  '

# Row 49
$ws.Range("E49").Value = '  -2.56%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0228'
$ws.Range("E50").Value = '  -3.19%  '

# Row 51
$ws.Range("E51").Value = '  -3.96%  '

